# Adds column S (year 2022) to the 11.5.1 "Number of deaths attributed to
# disasters" table, mirroring the layout already used by column R (2021).
#
# Style "recipes" (font / border / alignment) used by the existing sheet,
# keyed by the cellXfs index column R already uses on each row. We rebuild
# the same look on column S cell-by-cell because COM Copy/PasteSpecial on
# this host only transfers values, not per-cell formatting.
#   8  -> default font,            bottom border, no alignment override (row 3 rule)
#   12 -> bold Times New Roman 9,  bottom border, right/center align     (row 4 header)
#   13 -> bold Times New Roman 9,  theme text color, no border, right/center (row 5 totals)
#   20 -> plain Times New Roman 9, theme text color, no border, right/center
#   21 -> plain Times New Roman 9, theme text color, no border, right/center
#   24 -> plain Times New Roman 9, theme text color, bottom border, right/center (row 34 rule)
#   31 -> bold Times New Roman 9,  no theme color, no border, right/center
#   34 -> bold Times New Roman 9,  theme text color, no border, right/center

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Format-SCell($cell, $styleId) {

    switch ($styleId) {
        8 {
            $cell.Borders.Item(9).LineStyle = 1
            $cell.Borders.Item(9).Weight = -4138
        }
        12 {
            $cell.Font.Name = "Times New Roman"
            $cell.Font.Size = 9
            $cell.Font.Bold = $true
            $cell.HorizontalAlignment = -4152
            $cell.VerticalAlignment = -4108
            $cell.Borders.Item(9).LineStyle = 1
            $cell.Borders.Item(9).Weight = -4138
        }
        13 {
            $cell.Font.Name = "Times New Roman"
            $cell.Font.Size = 9
            $cell.Font.Bold = $true
            $cell.Font.ThemeColor = 1
            $cell.HorizontalAlignment = -4152
            $cell.VerticalAlignment = -4108
        }
        20 {
            $cell.Font.Name = "Times New Roman"
            $cell.Font.Size = 9
            $cell.Font.Bold = $false
            $cell.Font.ThemeColor = 1
            $cell.HorizontalAlignment = -4152
            $cell.VerticalAlignment = -4108
        }
        21 {
            $cell.Font.Name = "Times New Roman"
            $cell.Font.Size = 9
            $cell.Font.Bold = $false
            $cell.Font.ThemeColor = 1
            $cell.HorizontalAlignment = -4152
            $cell.VerticalAlignment = -4108
        }
        24 {
            $cell.Font.Name = "Times New Roman"
            $cell.Font.Size = 9
            $cell.Font.Bold = $false
            $cell.Font.ThemeColor = 1
            $cell.HorizontalAlignment = -4152
            $cell.VerticalAlignment = -4108
            $cell.Borders.Item(9).LineStyle = 1
            $cell.Borders.Item(9).Weight = -4138
        }
        31 {
            $cell.Font.Name = "Times New Roman"
            $cell.Font.Size = 9
            $cell.Font.Bold = $true
            $cell.HorizontalAlignment = -4152
            $cell.VerticalAlignment = -4108
        }
        34 {
            $cell.Font.Name = "Times New Roman"
            $cell.Font.Size = 9
            $cell.Font.Bold = $true
            $cell.Font.ThemeColor = 1
            $cell.HorizontalAlignment = -4152
            $cell.VerticalAlignment = -4108
        }
    }
}

# row -> style id, value ("-" marks the textual dash used for missing data)
$rows = @(
    @{ Row = 3;  Style = 8;  Value = $null },
    @{ Row = 4;  Style = 12; Value = 2022 },
    @{ Row = 5;  Style = 13; Value = 135 },
    @{ Row = 6;  Style = 20; Value = 99 },
    @{ Row = 7;  Style = 21; Value = 36 },
    @{ Row = 8;  Style = 31; Value = 97 },
    @{ Row = 9;  Style = 20; Value = 80 },
    @{ Row = 10; Style = 21; Value = 17 },
    @{ Row = 11; Style = 34; Value = 17 },
    @{ Row = 12; Style = 20; Value = 11 },
    @{ Row = 13; Style = 21; Value = 6 },
    @{ Row = 14; Style = 34; Value = 5 },
    @{ Row = 15; Style = 20; Value = 3 },
    @{ Row = 16; Style = 21; Value = 2 },
    @{ Row = 17; Style = 34; Value = "-" },
    @{ Row = 18; Style = 20; Value = "-" },
    @{ Row = 19; Style = 21; Value = "-" },
    @{ Row = 20; Style = 34; Value = 6 },
    @{ Row = 21; Style = 20; Value = 1 },
    @{ Row = 22; Style = 21; Value = 5 },
    @{ Row = 23; Style = 34; Value = "-" },
    @{ Row = 24; Style = 20; Value = "-" },
    @{ Row = 25; Style = 21; Value = "-" },
    @{ Row = 26; Style = 34; Value = 10 },
    @{ Row = 27; Style = 20; Value = 4 },
    @{ Row = 28; Style = 21; Value = 6 },
    @{ Row = 29; Style = 34; Value = "-" },
    @{ Row = 30; Style = 20; Value = "-" },
    @{ Row = 31; Style = 21; Value = "-" },
    @{ Row = 32; Style = 31; Value = "-" },
    @{ Row = 33; Style = 20; Value = "-" },
    @{ Row = 34; Style = 24; Value = "-" }
)

foreach ($entry in $rows) {
    $cell = $ws.Cells.Item($entry.Row, 19)   # column S = 19
    if ($null -ne $entry.Value) {
        $cell.Value = $entry.Value
    }
    Format-SCell $cell $entry.Style
}

# Matches the author's final selection after adding the 2022 column.
$ws.Range("S3").Select() | Out-Null
